$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append below the existing table (rows 29-33)
$rows = @(
    @{ A = "LIVE, TRAFFIC, METRO"; B = "13.737665552051753, 100.55985657894716"; C = "The Gaucho | Soi 19 | Sukhumvit Road"; D = "Bangkok"; E = "Thailand"; F = "xbBKbDwlR0E" },
    @{ A = "LIVE, TRAFFIC"; B = "16.07427386419858, 108.21753543894123"; C = "PTZ Camera Phuong Tran Da Nang Performance Costume"; D = "Da Nang"; E = "Vietnam"; F = "cB9Fs9UmcRU" },
    @{ A = "LIVE, TRAFFIC"; B = "16.074046511064708, 108.21727583355347"; C = "Camera at the back gate of Da Nang C Hospital"; D = "Da Nang"; E = "Vietnam"; F = "IXBTD4VgFF4" },
    @{ A = "LIVE, TRAFFIC"; B = "16.07413230110752, 108.21610230612718"; C = "Camera gate of Nguyen Hue school, Da Nang"; D = "Da Nang"; E = "Vietnam"; F = "Fu3nDsqC1J0" },
    @{ A = "LIVE, TRAFFIC, BUILDING"; B = "16.073735513035444, 108.21582517332482"; C = "Camera View of Danang Hospital Project"; D = "Da Nang"; E = "Vietnam"; F = "b6fkug3AmH4" }
)

$startRow = 29
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F

    # copy formatting from the row above so new rows match existing styling
    $ws.Range("A" + ($r - 1) + ":G" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":G" + $r).PasteSpecial(-4122)
}

# Fill the Status formula down through the new rows (last row left blank, per source)
$ws.Range("G29:G32").Formula = "=IsYouTubeVideoValid(F29)"

$excel.CutCopyMode = 0

$ws.Range("A34").Select()
